# Cash Flow workbook update:
#  - Switch the date column's custom number format from a date+time
#    format to a date-only format.
#  - Append four more transaction rows (rows 3-7) continuing the
#    running account-balance calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows (Date serial, Description, Amount, Account Balance)
$rows = @(
    @(44631, "Paycheck", 1350, 1350),
    @(44640, "Pizza",    -30,  1320),
    @(44645, "Paycheck", 1350, 2670),
    @(44659, "Paycheck", 1350, 4020),
    @(44673, "Paycheck", 1350, 5370)
)

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Keep the whole date column (including the new rows) on the same
# custom date format, now date-only instead of date+time.
$ws.Range("A2:A7").NumberFormat = "yyyy-mm-dd"
